$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C4 with new value and highlight it with a yellow fill
$ws.Range("C4").Value = 23800
$ws.Range("C4").Interior.Color = 65535

# Cascade formulas down C5:C9, each referencing the cell above
$ws.Range("C5").Formula = "=C4"
$ws.Range("C6").Formula = "=C5"
$ws.Range("C7").Formula = "=C6"
$ws.Range("C8").Formula = "=C7"
$ws.Range("C9").Formula = "=C8"

# Update the active selection to C4
$ws.Range("C4").Select()
